$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "Introduction`nid=`"Par1`">In the current study, we evaluated factors that increase the coronavirus disease (COVID-19) patient death rate by analyzing the data from two cohort hospitals.`n`n In addition, we studied whether underlying neurological diseases are risk factors for death.`n`n`nMethods`nid=`"Par2`">In this retrospective cohort study, we included 103 adult inpatients (aged " + [char]0x2265 + [char]0x2009 + "18" + [char]0x00A0 + "years).`n`n We evaluated differences in demographic data between surviving and non-surviving COVID-19 patients.`n`n`nResults`nid=`"Par3`">In a multivariate logistic analysis, age and the presence of chronic lung disease and Alzheimer" + [char]0x2019 + "s dementia (AD) were the only significant parameters for predicting COVID-19 non-survival (p" + [char]0x2009 + "&lt;" + [char]0x2009 + "0.05).`n`n However, hypertension, coronary vascular disease, dyslipidemia, chronic kidney disease, diabetes, and history of taking angiotensin II receptor blockers (ARBs) or angiotensin-converting enzyme (ACE) inhibitors, as well as nonsteroidal anti-inflammatory drugs (NSAIDs), were not significantly associated with the death of COVID-19 patients.`n`n The optimal cutoff value obtained from the maximum Youden index was 70 (sensitivity, 80.77%; specificity, 61.04%), and the odds ratio of non-survival increased 1.055 fold for every year of age.`n`n`nConclusions`nid=`"Par4`">Clinicians should closely monitor and manage the symptoms of COVID-19 patients who are over the age of 70" + [char]0x00A0 + "years or have chronic lung disease or AD.`n`n`n"

$ws.Range("E8").Value = "[Jong-moon%Hwang%NULL%0, Ju-Hyun%Kim%NULL%1, Jin-Sung%Park%NULL%1, Min Cheol%Chang%wheel633@ynu.ac.kr%1, Donghwi%Park%bdome@hanmail.net%1]"

$ws.Range("D9").Value = "id=`"Par1`">The natural history of coronavirus disease 2019 (COVID-19) has yet to be fully described.`n`n Here, we use patient-level data from the Information System for Research in Primary Care (SIDIAP) to summarise COVID-19 outcomes in Catalonia, Spain.`n`n We included 5,586,521 individuals from the general population.`n`n Of these, 102,002 had an outpatient diagnosis of COVID-19, 16,901 were hospitalised with COVID-19, and 5273 died after either being diagnosed or hospitalised with COVID-19 between 1st March and 6th May 2020. Older age, being male, and having comorbidities were all generally associated with worse outcomes.`n`n These findings demonstrate the continued need to protect those at high risk of poor outcomes, particularly older people, from COVID-19 and provide appropriate care for those who develop symptomatic disease.`n`n While risks of hospitalisation and death were lower for younger populations, there is a need to limit their role in community transmission.`n`n"

$ws.Range("E9").Value = "[Edward%Burn%NULL%1, Cristian%Teb" + [char]0x00E9 + "%NULL%2, Cristian%Teb" + [char]0x00E9 + "%NULL%0, Sergio%Fernandez-Bertolin%NULL%1, Maria%Aragon%NULL%1, Martina%Recalde%NULL%2, Martina%Recalde%NULL%0, Elena%Roel%NULL%1, Albert%Prats-Uribe%NULL%2, Albert%Prats-Uribe%NULL%0, Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%2, Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%0, Talita%Duarte-Salles%NULL%2, Talita%Duarte-Salles%NULL%0]"
